$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old rows 15-22 (years 2013-2020 leftover duplicates in the old layout)
$ws.Rows("15:22").Delete()

# Rewrite rows 2-14 with the shifted 2010-2022 data
$ws.Range("A2").Value = "2010年"
$ws.Range("B2").Value = 1122429.12188372
$ws.Range("C2").Value = 4318.2921685191
$ws.Range("D2").Value = 26808.3181027972

$ws.Range("A3").Value = "2011年"
$ws.Range("B3").Value = 1160716.81988533
$ws.Range("C3").Value = 4200.63315821866
$ws.Range("D3").Value = 26232.1640719654

$ws.Range("A4").Value = "2012年"
$ws.Range("B4").Value = 1244696.36849221
$ws.Range("C4").Value = 4219.29140686764
$ws.Range("D4").Value = 26606.2120737932

$ws.Range("A5").Value = "2013年"
$ws.Range("B5").Value = 1232371.60184321
$ws.Range("C5").Value = 4189.89992990652
$ws.Range("D5").Value = 26962.6818387152

$ws.Range("A6").Value = "2014年"
$ws.Range("B6").Value = 1202592.00951472
$ws.Range("C6").Value = 4200.40679494972
$ws.Range("D6").Value = 28051.3908781881

$ws.Range("A7").Value = "2015年"
$ws.Range("B7").Value = 1259132.00005221
$ws.Range("C7").Value = 4211.43921835381
$ws.Range("D7").Value = 28761.3978300345

$ws.Range("A8").Value = "2016年"
$ws.Range("B8").Value = 1319534.21879569
$ws.Range("C8").Value = 4264.95057582466
$ws.Range("D8").Value = 30005.3103718553

$ws.Range("A9").Value = "2017年"
$ws.Range("B9").Value = 1302190.57
$ws.Range("C9").Value = 4340.25988207897
$ws.Range("D9").Value = 30797.71

$ws.Range("A10").Value = "2018年"
$ws.Range("B10").Value = 1308936.04
$ws.Range("C10").Value = 4397.48323953924
$ws.Range("D10").Value = 31010.49

$ws.Range("A11").Value = "2019年"
$ws.Range("B11").Value = 1464062.2290628
$ws.Range("C11").Value = 4533.86955845411
$ws.Range("D11").Value = 31698.9018475957

$ws.Range("A12").Value = "2020年"
$ws.Range("B12").Value = 1557008.00414752
$ws.Range("C12").Value = 4565.45419939561
$ws.Range("D12").Value = 31941.2994664503

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = 4707.42935760008
$ws.Range("D13").Value = ""

$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = 4839.91497276742
$ws.Range("D14").Value = ""
